$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet immediately after "2021-Q3" (i.e.
#    right before "总计").
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $wsQ3)
$newSheet.Name = "2022-Q1"

# Copy header/body formatting (borders, bold, alignment) from the
# "2021-Q3" sheet so the new sheet matches its siblings.
$wsQ3.Range("A1:H2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "161030"
$newSheet.Range("C2").Value = "富国中证体育产业指数"
$newSheet.Range("D2").Value = "2.32"
$newSheet.Range("E2").Value = "93.75"
$newSheet.Range("F2").Value = "4.20"
$newSheet.Range("G2").Value = "0.0974"
$newSheet.Range("B2:G2").Style = "Normal"

$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" sheet: the date/count/value
#    columns shift down one row (oldest row 2020-Q4 becomes a brand-new
#    row 6), while column A stays the plain 0..4 running index.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$b2 = $wsTotal.Range("B2").Value()
$c2 = $wsTotal.Range("C2").Value()
$d2 = $wsTotal.Range("D2").Value()
$b3 = $wsTotal.Range("B3").Value()
$c3 = $wsTotal.Range("C3").Value()
$d3 = $wsTotal.Range("D3").Value()
$b4 = $wsTotal.Range("B4").Value()
$c4 = $wsTotal.Range("C4").Value()
$d4 = $wsTotal.Range("D4").Value()
$b5 = $wsTotal.Range("B5").Value()
$c5 = $wsTotal.Range("C5").Value()
$d5 = $wsTotal.Range("D5").Value()

# New row 6 (A6) picks up column-A's running-index formatting.
$wsTotal.Range("A5").Copy()
$wsTotal.Range("A6").PasteSpecial(-4122)
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = $b5
$wsTotal.Range("C6").Value = $c5
$wsTotal.Range("D6").Value = $d5

$wsTotal.Range("B5").Value = $b4
$wsTotal.Range("C5").Value = $c4
$wsTotal.Range("D5").Value = $d4

$wsTotal.Range("B4").Value = $b3
$wsTotal.Range("C4").Value = $c3
$wsTotal.Range("D4").Value = $d3

$wsTotal.Range("B3").Value = $b2
$wsTotal.Range("C3").Value = $c2
$wsTotal.Range("D3").Value = $d2

$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.1
